# Templates changes due to Director change
$d = $word.ActiveDocument

# 1) "Την" -> "Τη" (grammar fix) before the Φ.353.1/324/105657 decision reference.
$old1 = "Την με αριθ. Φ.353.1/324/105657"
$new1 = "Τη με αριθ. Φ.353.1/324/105657"
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $new1, 2)
Write-Output "Replace1: $found1"

# 2) Replace the whole Υπουργική Απόφαση reference (old Director decision -> new one).
#    This also folds the stray red-colored space run back into a single plain run,
#    matching the new text's single <w:t>.
$old2 = "Την με αριθ. Φ.350.2/1/32958/Ε3/27-2-2018  (ΑΔΑ:6Π414653ΠΣ-7ΕΝ) Υπουργική Απόφαση με θέμα: «Τοποθέτηση Περιφερειακών Διευθυντών Εκπαίδευσης»"
$new2 = "Τη με αριθ. Φ.351.1/11/48020/Ε3/28-3-2019 (ΑΔΑ: ΩΩΤΗ4653ΠΣ-ΒΔ3) Υπουργική Απόφαση με θέμα: «Τοποθέτηση Περιφερειακών Διευθυντών Εκπαίδευσης»"
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $new2, 2)
Write-Output "Replace2: $found2"
